$d = $word.ActiveDocument

# 1) The "_GoBack" bookmark currently sits at the end of the paragraph that
#    ends in "...mbios, con este podemos descartarlos y nuestros cambios se
#    iran." Remove it from there - it is going to be re-created at the end
#    of the new paragraph we add below (mirroring how Word keeps _GoBack
#    pinned to the most recent edit location).
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# 2) Append a new paragraph at the very end of the document documenting
#    "git push -u origin master".
$r = $d.Content
$r.Collapse(0)

$newParagraphXml = '<w:p><w:pPr><w:rPr><w:b/></w:rPr></w:pPr>' + `
'<w:proofErr w:type="spellStart"/>' + `
'<w:r><w:rPr><w:b/></w:rPr><w:t>git</w:t></w:r>' + `
'<w:proofErr w:type="spellEnd"/>' + `
'<w:r><w:rPr><w:b/></w:rPr><w:t xml:space="preserve"> </w:t></w:r>' + `
'<w:proofErr w:type="spellStart"/>' + `
'<w:r><w:rPr><w:b/></w:rPr><w:t>push</w:t></w:r>' + `
'<w:proofErr w:type="spellEnd"/>' + `
'<w:r><w:rPr><w:b/></w:rPr><w:t xml:space="preserve"> -u </w:t></w:r>' + `
'<w:proofErr w:type="spellStart"/>' + `
'<w:r><w:rPr><w:b/></w:rPr><w:t>origin</w:t></w:r>' + `
'<w:proofErr w:type="spellEnd"/>' + `
'<w:r><w:rPr><w:b/></w:rPr><w:t xml:space="preserve"> </w:t></w:r>' + `
'<w:proofErr w:type="gramStart"/>' + `
'<w:r><w:rPr><w:b/></w:rPr><w:t>master</w:t></w:r>' + `
'<w:proofErr w:type="gramEnd"/>' + `
'<w:r><w:rPr><w:b/></w:rPr><w:t xml:space="preserve">    #para de</w:t></w:r>' + `
'<w:r><w:rPr><w:b/></w:rPr><w:t>cirle el origen de nuestro proyecto</w:t></w:r>' + `
'<w:bookmarkStart w:id="0" w:name="_GoBack"/>' + `
'<w:bookmarkEnd w:id="0"/>' + `
'</w:p>'

$xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
'<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' + `
'<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' + `
$newParagraphXml + `
'</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$r.InsertXML($xml)
